# Cohort model rename: 'table' -> 'dataset'
$wb = $excel.ActiveWorkbook

$wsMolgenis = $wb.Worksheets.Item("molgenis")
$wsCollections = $wb.Worksheets.Item("MyCollections")
$wsDatasets = $wb.Worksheets.Item("MyTables")
$wsVariables = $wb.Worksheets.Item("MyVariables")
$wsHarmonisations = $wb.Worksheets.Item("MyVariableHarmonisations")

# 1) molgenis: "Tables"/"MyTables" row becomes "Datasets"/"MyDatasets"
$wsMolgenis.Range("B2").Value = "Datasets"
$wsMolgenis.Range("A2").Value = "MyDatasets"

# 2) MyVariables: "table" column header becomes "dataset"
$wsVariables.Range("B1").Value = "dataset"

# 3) MyVariableHarmonisations: "targetTable"/"sourceTable" become "targetDataset"/"sourceDataset"
$wsHarmonisations.Range("B1").Value = "targetDataset"
$wsHarmonisations.Range("E1").Value = "sourceDataset"

# 4) molgenis: "TableHarmonisations"/"MyTableHarmonisations" row becomes "DatasetHarmonisations"/"MyDatasetHarmonisations"
$wsMolgenis.Range("A10").Value = "MyDatasetHarmonisations"
$wsMolgenis.Range("B10").Value = "DatasetHarmonisations"

# 5) MyCollections: "name" column header becomes "acronym"
$wsCollections.Range("A1").Value = "acronym"

# 6) Rename the "MyTables" sheet to "MyDatasets" (keeps sheetId/rId, only the name changes;
#    sheet names are plain attributes, not shared strings, so ordering here is not significant)
$wsDatasets.Name = "MyDatasets"

# 7) Restore/update the cursor (selection) position on every sheet
$wsMolgenis.Activate()
$wsMolgenis.Range("B11").Select()

$wsCollections.Activate()
$wsCollections.Range("A2").Select()

$wsDatasets.Activate()
$wsDatasets.Range("B1").Select()

$wsVariables.Activate()
$wsVariables.Range("B2").Select()

$wsHarmonisations.Activate()
$wsHarmonisations.Range("E2").Select()

# 8) molgenis stays the active/visible tab, matching the original tabSelected="1"
$wsMolgenis.Activate()
